$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that get newly-assigned owner/status (previously blank B/C)
$newlyAssignedRows = @(181, 182, 183, 185, 186, 187, 188, 190, 192, 193)
foreach ($r in $newlyAssignedRows) {
    $ws.Cells.Item($r, 2).Value = "Erick Lim"
    $ws.Cells.Item($r, 3).Value = "Finished"
}

# Rows whose status moves from "In Progress" to "Finished"
$inProgressRows = @(194, 196, 197, 198, 199, 200, 201, 202, 203, 204)
foreach ($r in $inProgressRows) {
    $ws.Cells.Item($r, 3).Value = "Finished"
}

# Update the active selection to match the edited area
$ws.Range("C183").Select() | Out-Null
